# Update res_bus/vm_pu.xlsx sheet: case with 380 kV done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.041964649446291
$ws.Cells.Item(2, 4).Value = 1.044114952081765
$ws.Cells.Item(2, 5).Value = 1.049829942171636
$ws.Cells.Item(2, 6).Value = 1.060000298987649
$ws.Cells.Item(2, 9).Value = 1.040053609856913
$ws.Cells.Item(2, 10).Value = 1.047043020493194
$ws.Cells.Item(2, 11).Value = 1.046886911458517
$ws.Cells.Item(2, 12).Value = 1.052585910157114
$ws.Cells.Item(2, 13).Value = 1.0627282707336
$ws.Cells.Item(2, 14).Value = 1.019528603160958
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042995697018675
$ws.Cells.Item(3, 4).Value = 1.044901260075661
$ws.Cells.Item(3, 5).Value = 1.05082001121375
$ws.Cells.Item(3, 6).Value = 1.061240350381163
$ws.Cells.Item(3, 9).Value = 1.040318999402859
$ws.Cells.Item(3, 10).Value = 1.047719906683601
$ws.Cells.Item(3, 11).Value = 1.04748445678551
$ws.Cells.Item(3, 12).Value = 1.053387845190339
$ws.Cells.Item(3, 13).Value = 1.063781586190161
$ws.Cells.Item(3, 14).Value = 1.019757854908383
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043662932138472
$ws.Cells.Item(4, 4).Value = 1.045410049820293
$ws.Cells.Item(4, 5).Value = 1.051461119053009
$ws.Cells.Item(4, 6).Value = 1.06204370158996
$ws.Cells.Item(4, 9).Value = 1.040489561732674
$ws.Cells.Item(4, 10).Value = 1.048157385657466
$ws.Cells.Item(4, 11).Value = 1.047870448094591
$ws.Cells.Item(4, 12).Value = 1.05390660284762
$ws.Cells.Item(4, 13).Value = 1.06446354078635
$ws.Cells.Item(4, 14).Value = 1.019905915208565
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04394345665508
$ws.Cells.Item(5, 4).Value = 1.045623943476834
$ws.Cells.Item(5, 5).Value = 1.051730752394475
$ws.Cells.Item(5, 6).Value = 1.062381659371959
$ws.Cells.Item(5, 9).Value = 1.040560987676418
$ws.Cells.Item(5, 10).Value = 1.048341179137608
$ws.Cells.Item(5, 11).Value = 1.04803256042219
$ws.Cells.Item(5, 12).Value = 1.054124653050033
$ws.Cells.Item(5, 13).Value = 1.064750327791803
$ws.Cells.Item(5, 14).Value = 1.019968092351155
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.043990559065589
$ws.Cells.Item(6, 4).Value = 1.045659857055393
$ws.Cells.Item(6, 5).Value = 1.051776031558219
$ws.Cells.Item(6, 6).Value = 1.062438417404173
$ws.Cells.Item(6, 9).Value = 1.040572964072858
$ws.Cells.Item(6, 10).Value = 1.048372031672148
$ws.Cells.Item(6, 11).Value = 1.048059770496885
$ws.Cells.Item(6, 12).Value = 1.054161262552893
$ws.Cells.Item(6, 13).Value = 1.064798486067199
$ws.Cells.Item(6, 14).Value = 1.019978528209608
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.043666680446496
$ws.Cells.Item(7, 4).Value = 1.045412907884613
$ws.Cells.Item(7, 5).Value = 1.051464721468606
$ws.Cells.Item(7, 6).Value = 1.062048216498445
$ws.Cells.Item(7, 9).Value = 1.04049051722399
$ws.Cells.Item(7, 10).Value = 1.048159841997227
$ws.Cells.Item(7, 11).Value = 1.047872614870298
$ws.Cells.Item(7, 12).Value = 1.053909516584774
$ws.Cells.Item(7, 13).Value = 1.064467372481575
$ws.Cells.Item(7, 14).Value = 1.019906746287649
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.042313080630081
$ws.Cells.Item(8, 4).Value = 1.044380689058223
$ws.Cells.Item(8, 5).Value = 1.050164443839189
$ws.Cells.Item(8, 6).Value = 1.060419182591992
$ws.Cells.Item(8, 9).Value = 1.04014354008971
$ws.Cells.Item(8, 10).Value = 1.047271883177425
$ws.Cells.Item(8, 11).Value = 1.04708899127467
$ws.Cells.Item(8, 12).Value = 1.052856958425479
$ws.Cells.Item(8, 13).Value = 1.063084163155555
$ws.Cells.Item(8, 14).Value = 1.019606137932242
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039928460199267
$ws.Cells.Item(9, 4).Value = 1.042561780957705
$ws.Cells.Item(9, 5).Value = 1.047876772784714
$ws.Cells.Item(9, 6).Value = 1.057555917630924
$ws.Cells.Item(9, 9).Value = 1.039523226263905
$ws.Cells.Item(9, 10).Value = 1.045703275695203
$ws.Cells.Item(9, 11).Value = 1.045703100899467
$ws.Cells.Item(9, 12).Value = 1.051001087612939
$ws.Cells.Item(9, 13).Value = 1.060649744282764
$ws.Cells.Item(9, 14).Value = 1.019074282986955
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.038339103645887
$ws.Cells.Item(10, 4).Value = 1.04134920265555
$ws.Cells.Item(10, 5).Value = 1.046354076622426
$ws.Cells.Item(10, 6).Value = 1.055651961266321
$ws.Cells.Item(10, 9).Value = 1.03910371140401
$ws.Cells.Item(10, 10).Value = 1.044654916943254
$ws.Cells.Item(10, 11).Value = 1.044775797400694
$ws.Cells.Item(10, 12).Value = 1.049763084667815
$ws.Cells.Item(10, 13).Value = 1.059028784428888
$ws.Cells.Item(10, 14).Value = 1.01871827962548
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037650984193456
$ws.Cells.Item(11, 4).Value = 1.040824155336858
$ws.Cells.Item(11, 5).Value = 1.045695306707466
$ws.Cells.Item(11, 6).Value = 1.054828680504469
$ws.Cells.Item(11, 9).Value = 1.038920641775578
$ws.Cells.Item(11, 10).Value = 1.044200345685051
$ws.Cells.Item(11, 11).Value = 1.044373466596
$ws.Cells.Item(11, 12).Value = 1.049226836388086
$ws.Cells.Item(11, 13).Value = 1.058327358314479
$ws.Cells.Item(11, 14).Value = 1.018563787830623
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.037395397995875
$ws.Cells.Item(12, 4).Value = 1.04062913086269
$ws.Cells.Item(12, 5).Value = 1.045450695469138
$ws.Cells.Item(12, 6).Value = 1.054523048784257
$ws.Cells.Item(12, 9).Value = 1.038852428751035
$ws.Cells.Item(12, 10).Value = 1.044031403787536
$ws.Cells.Item(12, 11).Value = 1.04422390263003
$ws.Cells.Item(12, 12).Value = 1.049027622027275
$ws.Cells.Item(12, 13).Value = 1.058066886146735
$ws.Cells.Item(12, 14).Value = 1.018506351716425
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.037450221565615
$ws.Cells.Item(13, 4).Value = 1.040670964200936
$ws.Cells.Item(13, 5).Value = 1.045503161553317
$ws.Cells.Item(13, 6).Value = 1.05458860007857
$ws.Cells.Item(13, 9).Value = 1.038867070305566
$ws.Cells.Item(13, 10).Value = 1.044067646646746
$ws.Cells.Item(13, 11).Value = 1.04425599005413
$ws.Cells.Item(13, 12).Value = 1.049070355446889
$ws.Cells.Item(13, 13).Value = 1.058122755207754
$ws.Cells.Item(13, 14).Value = 1.01851867426342
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037629857115894
$ws.Cells.Item(14, 4).Value = 1.040808034520455
$ws.Cells.Item(14, 5).Value = 1.045675085327507
$ws.Cells.Item(14, 6).Value = 1.054803413397838
$ws.Cells.Item(14, 9).Value = 1.038915007606643
$ws.Cells.Item(14, 10).Value = 1.044186382810433
$ws.Cells.Item(14, 11).Value = 1.044361106047045
$ws.Cells.Item(14, 12).Value = 1.049210369836076
$ws.Cells.Item(14, 13).Value = 1.058305826204329
$ws.Cells.Item(14, 14).Value = 1.018559041182982
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037740538086768
$ws.Cells.Item(15, 4).Value = 1.040892488276355
$ws.Cells.Item(15, 5).Value = 1.045781024533187
$ws.Cells.Item(15, 6).Value = 1.054935789669427
$ws.Cells.Item(15, 9).Value = 1.03894451516207
$ws.Cells.Item(15, 10).Value = 1.044259527638244
$ws.Cells.Item(15, 11).Value = 1.044425855530755
$ws.Cells.Item(15, 12).Value = 1.049296633625767
$ws.Cells.Item(15, 13).Value = 1.058418631381357
$ws.Cells.Item(15, 14).Value = 1.018583905821573
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038384773571932
$ws.Cells.Item(16, 4).Value = 1.041384048501481
$ws.Cells.Item(16, 5).Value = 1.046397808967434
$ws.Cells.Item(16, 6).Value = 1.055706623769811
$ws.Cells.Item(16, 9).Value = 1.039115831284722
$ws.Cells.Item(16, 10).Value = 1.044685072173072
$ws.Cells.Item(16, 11).Value = 1.044802481880234
$ws.Cells.Item(16, 12).Value = 1.04979866979687
$ws.Cells.Item(16, 13).Value = 1.059075345430145
$ws.Cells.Item(16, 14).Value = 1.018728525582166
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038788907270975
$ws.Cells.Item(17, 4).Value = 1.041692393433321
$ws.Cells.Item(17, 5).Value = 1.046784853480305
$ws.Cells.Item(17, 6).Value = 1.05619045373685
$ws.Cells.Item(17, 9).Value = 1.03922291405559
$ws.Cells.Item(17, 10).Value = 1.044951837691174
$ws.Cells.Item(17, 11).Value = 1.045038514947573
$ws.Cells.Item(17, 12).Value = 1.050113534263988
$ws.Cells.Item(17, 13).Value = 1.05948740771543
$ws.Cells.Item(17, 14).Value = 1.018819150780644
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.039024639598051
$ws.Cells.Item(18, 4).Value = 1.041872246471539
$ws.Cells.Item(18, 5).Value = 1.047010664760845
$ws.Cells.Item(18, 6).Value = 1.056472774077811
$ws.Cells.Item(18, 9).Value = 1.039285236906988
$ws.Cells.Item(18, 10).Value = 1.045107377200559
$ws.Cells.Item(18, 11).Value = 1.04517611162545
$ws.Cells.Item(18, 12).Value = 1.050297171553551
$ws.Cells.Item(18, 13).Value = 1.059727801105241
$ws.Cells.Item(18, 14).Value = 1.018871978115073
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.039105019594422
$ws.Cells.Item(19, 4).Value = 1.041933571808488
$ws.Cells.Item(19, 5).Value = 1.047087669862
$ws.Cells.Item(19, 6).Value = 1.056569056780419
$ws.Cells.Item(19, 9).Value = 1.039306464197236
$ws.Cells.Item(19, 10).Value = 1.045160401924963
$ws.Cells.Item(19, 11).Value = 1.045223015401262
$ws.Cells.Item(19, 12).Value = 1.050359784103095
$ws.Cells.Item(19, 13).Value = 1.059809776587756
$ws.Cells.Item(19, 14).Value = 1.018889985302737
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.038745546727164
$ws.Cells.Item(20, 4).Value = 1.041659310869341
$ws.Cells.Item(20, 5).Value = 1.046743321602599
$ws.Cells.Item(20, 6).Value = 1.056138531934226
$ws.Cells.Item(20, 9).Value = 1.039211439228492
$ws.Cells.Item(20, 10).Value = 1.044923222513116
$ws.Cells.Item(20, 11).Value = 1.045013198839748
$ws.Cells.Item(20, 12).Value = 1.050079754131098
$ws.Cells.Item(20, 13).Value = 1.059443192739182
$ws.Cells.Item(20, 14).Value = 1.018809430947542
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037576958595237
$ws.Cells.Item(21, 4).Value = 1.04076767068364
$ws.Cells.Item(21, 5).Value = 1.045624455708302
$ws.Cells.Item(21, 6).Value = 1.054740151502385
$ws.Cells.Item(21, 9).Value = 1.038900897140451
$ws.Cells.Item(21, 10).Value = 1.044151420561476
$ws.Cells.Item(21, 11).Value = 1.044330155329334
$ws.Cells.Item(21, 12).Value = 1.049169139858489
$ws.Cells.Item(21, 13).Value = 1.058251914465028
$ws.Cells.Item(21, 14).Value = 1.018547155537889
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036842290131787
$ws.Cells.Item(22, 4).Value = 1.040207069863343
$ws.Cells.Item(22, 5).Value = 1.044921473470077
$ws.Cells.Item(22, 6).Value = 1.053861925663554
$ws.Cells.Item(22, 9).Value = 1.03870441576468
$ws.Cells.Item(22, 10).Value = 1.043665614813611
$ws.Cells.Item(22, 11).Value = 1.043900002287685
$ws.Cells.Item(22, 12).Value = 1.048596439207215
$ws.Cells.Item(22, 13).Value = 1.057503307933668
$ws.Cells.Item(22, 14).Value = 1.018381957654982
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.037231745343299
$ws.Cells.Item(23, 4).Value = 1.040504253997738
$ws.Cells.Item(23, 5).Value = 1.045294090919868
$ws.Cells.Item(23, 6).Value = 1.05432739601887
$ws.Cells.Item(23, 9).Value = 1.038808690945094
$ws.Cells.Item(23, 10).Value = 1.043923201121423
$ws.Cells.Item(23, 11).Value = 1.044128100524587
$ws.Cells.Item(23, 12).Value = 1.048900053983368
$ws.Cells.Item(23, 13).Value = 1.057900120889484
$ws.Cells.Item(23, 14).Value = 1.018469560101394
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038765139465626
$ws.Cells.Item(24, 4).Value = 1.041674259453913
$ws.Cells.Item(24, 5).Value = 1.04676208790107
$ws.Cells.Item(24, 6).Value = 1.056161992820181
$ws.Cells.Item(24, 9).Value = 1.039216624631313
$ws.Cells.Item(24, 10).Value = 1.044936152666578
$ws.Cells.Item(24, 11).Value = 1.045024638338285
$ws.Cells.Item(24, 12).Value = 1.05009501797405
$ws.Cells.Item(24, 13).Value = 1.059463171446242
$ws.Cells.Item(24, 14).Value = 1.018813823023093
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040544872213097
$ws.Cells.Item(25, 4).Value = 1.043032009831304
$ws.Cells.Item(25, 5).Value = 1.048467765012474
$ws.Cells.Item(25, 6).Value = 1.058295277317452
$ws.Cells.Item(25, 9).Value = 1.039684645560286
$ws.Cells.Item(25, 10).Value = 1.046109260847172
$ws.Cells.Item(25, 11).Value = 1.046061983078429
$ws.Cells.Item(25, 12).Value = 1.051481008149697
$ws.Cells.Item(25, 13).Value = 1.061278748945669
$ws.Cells.Item(25, 14).Value = 1.019212033507536
